# Auto-generated Excel COM-interop script
# Applies numeric refresh updates to Leve profit-calculation columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled-runner diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 239.6
$ws.Range("I4").Value = 239.6
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 239.6
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -125.6
$ws.Range("N4").ClearContents()
$ws.Range("H41").Value = 666.9167
$ws.Range("I41").Value = 863.1429000000001
$ws.Range("J41").Value = 392.2
$ws.Range("K41").Value = 863.1429000000001
$ws.Range("L41").Value = 392.2
$ws.Range("M41").Value = -423.1429000000001
$ws.Range("N41").Value = -1272.2
$ws.Range("H74").Value = 5154.2
$ws.Range("I74").Value = 5154.2
$ws.Range("K74").Value = 5154.2
$ws.Range("M74").Value = -4218.2
$ws.Range("H77").Value = 5154.2
$ws.Range("I77").Value = 5154.2
$ws.Range("K77").Value = 25771
$ws.Range("M77").Value = -21091
$ws.Range("H80").Value = 2624.4546
$ws.Range("I80").Value = 1540.6666
$ws.Range("K80").Value = 4621.9998
$ws.Range("M80").Value = -3623.9998
$ws.Range("H83").Value = 2624.4546
$ws.Range("I83").Value = 1540.6666
$ws.Range("K83").Value = 13865.9994
$ws.Range("M83").Value = -8873.999400000001
$ws.Range("H107").Value = 1549.6666
$ws.Range("I107").Value = 1549.6666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1549.6666
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 370.3334
$ws.Range("N107").ClearContents()
$ws.Range("H112").Value = 1918.6666
$ws.Range("J112").Value = 1838.5
$ws.Range("L112").Value = 5515.5
$ws.Range("N112").Value = -7731.5
$ws.Range("H137").Value = 2314.2
$ws.Range("I137").Value = 2149
$ws.Range("J137").Value = 2424.3333
$ws.Range("K137").Value = 6447
$ws.Range("L137").Value = 7272.999899999999
$ws.Range("M137").Value = -3897
$ws.Range("N137").Value = -12372.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value = 2994.4138
$ws.Range("I32").Value = 3186.1924
$ws.Range("K32").Value = 3186.1924
$ws.Range("M32").Value = -2899.1924
$ws.Range("H37").Value = 2500
$ws.Range("J37").Value = 2500
$ws.Range("L37").Value = 2500
$ws.Range("N37").Value = -3046
$ws.Range("H44").Value = 27032.5
$ws.Range("J44").Value = 27032.5
$ws.Range("L44").Value = 27032.5
$ws.Range("N44").Value = -28008.5
$ws.Range("H55").Value = 23173
$ws.Range("I55").Value = 19999.5
$ws.Range("J55").Value = 24442.4
$ws.Range("K55").Value = 19999.5
$ws.Range("L55").Value = 24442.4
$ws.Range("M55").Value = -19684.5
$ws.Range("N55").Value = -25072.4
$ws.Range("H61").Value = 3433.4707
$ws.Range("I61").Value = 3137.5
$ws.Range("K61").Value = 3137.5
$ws.Range("M61").Value = -2925.5
$ws.Range("H74").Value = 2106478.2
$ws.Range("J74").Value = 5557556
$ws.Range("L74").Value = 5557556
$ws.Range("N74").Value = -5559304
$ws.Range("H77").Value = 2106478.2
$ws.Range("J77").Value = 5557556
$ws.Range("L77").Value = 27787780
$ws.Range("N77").Value = -27796516
$ws.Range("H102").Value = 1178
$ws.Range("I102").Value = 1097.5
$ws.Range("K102").Value = 1097.5
$ws.Range("M102").Value = 524.5
$ws.Range("H132").Value = 12823008
$ws.Range("I132").Value = 2577.5938
$ws.Range("K132").Value = 7732.7814
$ws.Range("M132").Value = -5202.7814
$ws.Range("H136").Value = 3433.4707
$ws.Range("I136").Value = 3137.5
$ws.Range("K136").Value = 9412.5
$ws.Range("M136").Value = -6862.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3751.1892
$ws.Range("I86").Value = 1584.9524
$ws.Range("K86").Value = 1584.9524
$ws.Range("M86").Value = -461.9523999999999
$ws.Range("H89").Value = 3751.1892
$ws.Range("I89").Value = 1584.9524
$ws.Range("K89").Value = 7924.762
$ws.Range("M89").Value = -2308.762
$ws.Range("H94").Value = 1330.625
$ws.Range("J94").Value = 579.6667
$ws.Range("L94").Value = 579.6667
$ws.Range("N94").Value = -1481.6667
$ws.Range("H105").Value = 3635.0908
$ws.Range("I105").Value = 1765
$ws.Range("K105").Value = 1765
$ws.Range("M105").Value = -18

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2776.1724
$ws.Range("I31").Value = 2616.7058
$ws.Range("J31").Value = 3002.0833
$ws.Range("K31").Value = 2616.7058
$ws.Range("L31").Value = 3002.0833
$ws.Range("M31").Value = -2321.7058
$ws.Range("N31").Value = -3592.0833
$ws.Range("H34").Value = 2776.1724
$ws.Range("I34").Value = 2616.7058
$ws.Range("J34").Value = 3002.0833
$ws.Range("K34").Value = 2616.7058
$ws.Range("L34").Value = 3002.0833
$ws.Range("M34").Value = -2414.7058
$ws.Range("N34").Value = -3406.0833
$ws.Range("H62").Value = 82925.5
$ws.Range("I62").Value = 3902.5
$ws.Range("J62").Value = 109266.5
$ws.Range("K62").Value = 3902.5
$ws.Range("L62").Value = 109266.5
$ws.Range("M62").Value = -3278.5
$ws.Range("N62").Value = -110514.5
$ws.Range("H65").Value = 82925.5
$ws.Range("I65").Value = 3902.5
$ws.Range("J65").Value = 109266.5
$ws.Range("K65").Value = 19512.5
$ws.Range("L65").Value = 546332.5
$ws.Range("M65").Value = -16392.5
$ws.Range("N65").Value = -552572.5
$ws.Range("H99").Value = 2324.2144
$ws.Range("I99").Value = 1878.6666
$ws.Range("J99").Value = 2658.375
$ws.Range("K99").Value = 1878.6666
$ws.Range("L99").Value = 2658.375
$ws.Range("M99").Value = -380.6666
$ws.Range("N99").Value = -5654.375
$ws.Range("H126").Value = 2324.2144
$ws.Range("I126").Value = 1878.6666
$ws.Range("J126").Value = 2658.375
$ws.Range("K126").Value = 5635.9998
$ws.Range("L126").Value = 7975.125
$ws.Range("M126").Value = -3165.9998
$ws.Range("N126").Value = -12915.125
$ws.Range("H132").Value = 10526.448
$ws.Range("J132").Value = 8995.5
$ws.Range("L132").Value = 26986.5
$ws.Range("N132").Value = -32046.5
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 971
$ws.Range("I8").Value = 971
$ws.Range("K8").Value = 2913
$ws.Range("M8").Value = -2774
$ws.Range("H10").Value = 297.83334
$ws.Range("I10").Value = 317
$ws.Range("K10").Value = 951
$ws.Range("M10").Value = -812
$ws.Range("H103").Value = 637.2308
$ws.Range("I103").Value = 326.42856
$ws.Range("K103").Value = 979.28568
$ws.Range("M103").Value = -100.28568

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16592.6
$ws.Range("I80").Value = 5988
$ws.Range("K80").Value = 5988
$ws.Range("M80").Value = -4990
$ws.Range("H83").Value = 16592.6
$ws.Range("I83").Value = 5988
$ws.Range("K83").Value = 29940
$ws.Range("M83").Value = -24948
$ws.Range("H103").Value = 69998.336
$ws.Range("J103").Value = 69998.336
$ws.Range("L103").Value = 69998.336
$ws.Range("N103").Value = -72342.336
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 69998.5
$ws.Range("J135").Value = 69998.5
$ws.Range("L135").Value = 69998.5
$ws.Range("N135").Value = -80138.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1736.3889
$ws.Range("I16").Value = 1683.7333
$ws.Range("J16").Value = 1999.6666
$ws.Range("K16").Value = 1683.7333
$ws.Range("L16").Value = 1999.6666
$ws.Range("M16").Value = -1513.7333
$ws.Range("N16").Value = -2339.6666
$ws.Range("H82").Value = 685.8889
$ws.Range("I82").Value = 646.4
$ws.Range("J82").Value = 883.3333
$ws.Range("K82").Value = 646.4
$ws.Range("L82").Value = 883.3333
$ws.Range("M82").Value = -285.4
$ws.Range("N82").Value = -1605.3333
$ws.Range("H85").Value = 685.8889
$ws.Range("I85").Value = 646.4
$ws.Range("J85").Value = 883.3333
$ws.Range("K85").Value = 646.4
$ws.Range("L85").Value = 883.3333
$ws.Range("M85").Value = 601.6
$ws.Range("N85").Value = -3379.3333
$ws.Range("H136").Value = 22225134
$ws.Range("I136").Value = 2559.842
$ws.Range("K136").Value = 7679.526
$ws.Range("M136").Value = -5129.526

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 29998.334
$ws.Range("I34").Value = 29998.334
$ws.Range("K34").Value = 29998.334
$ws.Range("M34").Value = -29795.334
$ws.Range("H100").Value = 1583.96
$ws.Range("I100").Value = 1638.75
$ws.Range("J100").Value = 1364.8
$ws.Range("K100").Value = 3277.5
$ws.Range("L100").Value = 2729.6
$ws.Range("M100").Value = -2736.5
$ws.Range("N100").Value = -3811.6
$ws.Range("H132").Value = 262
$ws.Range("I132").Value = 244.5
$ws.Range("K132").Value = 733.5
$ws.Range("M132").Value = 1796.5
$ws.Range("H136").Value = 1655.5714
$ws.Range("I136").Value = 1473.1666
$ws.Range("K136").Value = 4419.4998
$ws.Range("M136").Value = -1869.4998
